$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data rows 2-10 (Q0..Q8) with new computed values ---

# Row 2 (Q0)
$ws.Range("B2").Value = 0.4257792306779735
$ws.Range("C2").Value = 2.217877629047371
$ws.Range("D2").Value = 14.65900008690008
$ws.Range("E2").Value = 3.828707365012384
$ws.Range("F2").Value = 3.842820107199375
$ws.Range("G2").Value = 51

# Row 3 (Q1)
$ws.Range("B3").Value = 0.7099789703697202
$ws.Range("C3").Value = 2.408638386592436
$ws.Range("D3").Value = 14.21756035022525
$ws.Range("E3").Value = 3.770618032925802
$ws.Range("F3").Value = 3.740769621219512
$ws.Range("G3").Value = 50

# Row 4 (Q2)
$ws.Range("B4").Value = 0.5489223891954144
$ws.Range("C4").Value = 2.142174787177605
$ws.Range("D4").Value = 13.77538404339546
$ws.Range("E4").Value = 3.711520449006776
$ws.Range("F4").Value = 3.708743454234583
$ws.Range("G4").Value = 49

# Row 5 (Q3)
$ws.Range("B5").Value = 0.7318547725335144
$ws.Range("C5").Value = 2.680487997555033
$ws.Range("D5").Value = 18.5702968988901
$ws.Range("E5").Value = 4.309326733828627
$ws.Range("F5").Value = 4.291666599746589
$ws.Range("G5").Value = 48

# Row 6 (Q4)
$ws.Range("B6").Value = 0.6680602970490725
$ws.Range("C6").Value = 2.486664053058593
$ws.Range("D6").Value = 14.69638859368018
$ws.Range("E6").Value = 3.833586909629176
$ws.Range("F6").Value = 3.815739570344253
$ws.Range("G6").Value = 47

# Row 7 (Q5)
$ws.Range("B7").Value = 0.6768152876208176
$ws.Range("C7").Value = 2.661179530454853
$ws.Range("D7").Value = 17.12325255943136
$ws.Range("E7").Value = 4.138025200434545
$ws.Range("F7").Value = 4.127409698292647
$ws.Range("G7").Value = 46

# Row 8 (Q6)
$ws.Range("B8").Value = 0.5558206604367482
$ws.Range("C8").Value = 2.164092278013886
$ws.Range("D8").Value = 13.2712682877262
$ws.Range("E8").Value = 3.642975197242798
$ws.Range("F8").Value = 3.641006746554708
$ws.Range("G8").Value = 45

# Row 9 (Q7)
$ws.Range("B9").Value = 0.6185392898960298
$ws.Range("C9").Value = 2.459293877902776
$ws.Range("D9").Value = 13.75218320165127
$ws.Range("E9").Value = 3.70839361471396
$ws.Range("F9").Value = 3.69871776441472
$ws.Range("G9").Value = 44

# Row 10 (Q8) - note: F10 did not previously exist and is newly added here
$ws.Range("B10").Value = 0.7423495001570053
$ws.Range("C10").Value = 2.530096135466969
$ws.Range("D10").Value = 17.10393720569436
$ws.Range("E10").Value = 4.135690656431446
$ws.Range("F10").Value = 4.116669771413056
$ws.Range("G10").Value = 43

# --- Append new row 11 (Q9) ---
$ws.Range("A11").Value = "Q9"
$ws.Range("B11").Value = 0.5511858627957305
$ws.Range("C11").Value = 2.216421366742779
$ws.Range("D11").Value = 14.11284219713519
$ws.Range("E11").Value = 3.756706296363237
$ws.Range("F11").Value = 3.761095865066883
$ws.Range("G11").Value = 42

# Match the formatting used by the other label cells in column A
# (bold font, thin box border, centered horizontal/top vertical alignment)
# by copying the format from the cell directly above it.
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122)
